$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 346.33334
$ws.Range("I28").Value = 346.33334
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 346.33334
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = $null
$ws.Range("N28").Value = 138.66666
$ws.Range("H129").Value = 4387089.5
$ws.Range("J129").Value = 1068.4445
$ws.Range("L129").Value = 3205.3335
$ws.Range("N129").Value = -13205.3335
$ws.Range("H137").Value = 3617.1333
$ws.Range("I137").Value = 4335.5557
$ws.Range("J137").Value = 2539.5
$ws.Range("K137").Value = 13006.6671
$ws.Range("L137").Value = 7618.5
$ws.Range("M137").Value = -10456.6671
$ws.Range("N137").Value = -12718.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3290.31
$ws.Range("I32").Value = 3290.31
$ws.Range("K32").Value = 3290.31
$ws.Range("M32").Value = -3003.31
$ws.Range("H61").Value = 4375.5
$ws.Range("I61").Value = 3739.3333
$ws.Range("J61").Value = 5011.6665
$ws.Range("K61").Value = 3739.3333
$ws.Range("L61").Value = 5011.6665
$ws.Range("M61").Value = -3527.3333
$ws.Range("N61").Value = -5435.6665
$ws.Range("H103").Value = 25185.5
$ws.Range("J103").Value = 25185.5
$ws.Range("L103").Value = 25185.5
$ws.Range("N103").Value = -27529.5
$ws.Range("H132").Value = 2095.111
$ws.Range("I132").Value = 1555.6904
$ws.Range("K132").Value = 4667.0712
$ws.Range("M132").Value = -2137.0712
$ws.Range("H133").Value = 25720.334
$ws.Range("J133").Value = 25720.334
$ws.Range("L133").Value = 25720.334
$ws.Range("N133").Value = -30780.334
$ws.Range("H134").Value = 29880
$ws.Range("J134").Value = 29880
$ws.Range("L134").Value = 29880
$ws.Range("N134").Value = -40020
$ws.Range("H135").Value = 32330
$ws.Range("J135").Value = 32330
$ws.Range("L135").Value = 32330
$ws.Range("N135").Value = -42470
$ws.Range("H136").Value = 4375.5
$ws.Range("I136").Value = 3739.3333
$ws.Range("J136").Value = 5011.6665
$ws.Range("K136").Value = 11217.9999
$ws.Range("L136").Value = 15034.9995
$ws.Range("M136").Value = -8667.999899999999
$ws.Range("N136").Value = -20134.9995
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = $null
$ws.Range("N137").Value = 0
$ws.Range("H138").Value = 98214.5
$ws.Range("J138").Value = 98214.5
$ws.Range("L138").Value = 98214.5
$ws.Range("N138").Value = -108494.5
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = $null
$ws.Range("N139").Value = 0
$ws.Range("H140").Value = 29700
$ws.Range("J140").Value = 29700
$ws.Range("L140").Value = 29700
$ws.Range("N140").Value = -40060
$ws.Range("H141").Value = 29583.334
$ws.Range("J141").Value = 29583.334
$ws.Range("L141").Value = 29583.334
$ws.Range("N141").Value = -39943.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 32450
$ws.Range("J103").Value = 32450
$ws.Range("L103").Value = 32450
$ws.Range("N103").Value = -34794
$ws.Range("H107").Value = 2169.9473
$ws.Range("I107").Value = 1849.25
$ws.Range("J107").Value = 2719.7144
$ws.Range("K107").Value = 1849.25
$ws.Range("L107").Value = 2719.7144
$ws.Range("M107").Value = 70.75
$ws.Range("N107").Value = -6559.7144
$ws.Range("H134").Value = 3660.1277
$ws.Range("I134").Value = 3378.3555
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 10135.0665
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -7600.066500000001
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2331.353
$ws.Range("I122").Value = 1939.3704
$ws.Range("K122").Value = 5818.1112
$ws.Range("M122").Value = -3368.1112
$ws.Range("H132").Value = 2339.2727
$ws.Range("I132").Value = 1868.8966
$ws.Range("J132").Value = 5749.5
$ws.Range("K132").Value = 5606.6898
$ws.Range("L132").Value = 17248.5
$ws.Range("M132").Value = -3076.6898
$ws.Range("N132").Value = -22308.5
$ws.Range("H134").Value = 11366066
$ws.Range("I134").Value = 14707639
$ws.Range("J134").Value = 4719.9
$ws.Range("K134").Value = 44122917
$ws.Range("L134").Value = 14159.7
$ws.Range("M134").Value = -44120382
$ws.Range("N134").Value = -19229.7
$ws.Range("H141").Value = 25638.889
$ws.Range("J141").Value = 25638.889
$ws.Range("L141").Value = 25638.889
$ws.Range("N141").Value = -35998.889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 748.0238000000001
$ws.Range("I5").Value = 446.58975
$ws.Range("K5").Value = 1339.76925
$ws.Range("M5").Value = -1227.76925
$ws.Range("H97").Value = 2037.5
$ws.Range("J97").Value = 3000
$ws.Range("L97").Value = 9000
$ws.Range("N97").Value = -9992
$ws.Range("H107").Value = 789.4666999999999
$ws.Range("I107").Value = 474.125
$ws.Range("J107").Value = 1149.8572
$ws.Range("K107").Value = 1422.375
$ws.Range("L107").Value = 3449.5716
$ws.Range("M107").Value = 497.625
$ws.Range("N107").Value = -7289.571599999999
$ws.Range("H131").Value = 1222.5
$ws.Range("J131").Value = 1163.2927
$ws.Range("L131").Value = 3489.8781
$ws.Range("N131").Value = -13569.8781
$ws.Range("H135").Value = 748.0238000000001
$ws.Range("I135").Value = 446.58975
$ws.Range("K135").Value = 4019.30775
$ws.Range("M135").Value = -1484.30775
$ws.Range("H136").Value = 2227.3572
$ws.Range("J136").Value = 3394.6
$ws.Range("L136").Value = 10183.8
$ws.Range("N136").Value = -20383.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2345.5
$ws.Range("I113").Value = 884.8333
$ws.Range("J113").Value = 3806.1667
$ws.Range("K113").Value = 884.8333
$ws.Range("L113").Value = 3806.1667
$ws.Range("M113").Value = 1285.1667
$ws.Range("N113").Value = -8146.1667
$ws.Range("H126").Value = 558700.25
$ws.Range("I126").Value = 1976.7778
$ws.Range("J126").Value = 1115423.8
$ws.Range("K126").Value = 5930.3334
$ws.Range("L126").Value = 3346271.4
$ws.Range("M126").Value = -3460.3334
$ws.Range("N126").Value = -3351211.4
$ws.Range("H132").Value = 2431.4143
$ws.Range("I132").Value = 2019.0754
$ws.Range("J132").Value = 3716.9412
$ws.Range("K132").Value = 6057.2262
$ws.Range("L132").Value = 11150.8236
$ws.Range("M132").Value = -3527.2262
$ws.Range("N132").Value = -16210.8236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2124.5894
$ws.Range("I132").Value = 1419.6487
$ws.Range("J132").Value = 3497.3684
$ws.Range("K132").Value = 4258.9461
$ws.Range("L132").Value = 10492.1052
$ws.Range("M132").Value = -1728.9461
$ws.Range("N132").Value = -15552.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11400.0625
$ws.Range("I132").Value = 2766.6667
$ws.Range("J132").Value = 71833.836
$ws.Range("K132").Value = 8300.000100000001
$ws.Range("L132").Value = 215501.508
$ws.Range("M132").Value = -5770.000100000001
$ws.Range("N132").Value = -220561.508
